# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1163
$ws.Range("F8").Value = 1096
$ws.Range("G8").Value = 79.2
$ws.Range("F9").Value = 1657
$ws.Range("F10").Value = 6156
$ws.Range("F11").Value = 120
$ws.Range("F16").Value = 6406
$ws.Range("F21").Value = 1682
$ws.Range("F23").Value = 6
$ws.Range("F26").Value = 1484
$ws.Range("F33").Value = 3881

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 323
$ws.Range("F8").Value = 402

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9519
$ws.Range("F3").Value = 2245
$ws.Range("F4").Value = 648
$ws.Range("F5").Value = 218

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9519
$ws.Range("F3").Value = 2245
$ws.Range("F4").Value = 648
$ws.Range("F7").Value = 1163
$ws.Range("F11").Value = 323
$ws.Range("F12").Value = 218
$ws.Range("F13").Value = 1657
$ws.Range("F14").Value = 6156
$ws.Range("F15").Value = 120
$ws.Range("F23").Value = 6406
$ws.Range("F28").Value = 1682
$ws.Range("F32").Value = 1484
$ws.Range("F44").Value = 3881
